# Auto-generated edit script: refreshes scraped market-price columns
# (H..N) across several leve rows on multiple sheets, per the scheduled
# market-data runner. Generated from the authoritative cell diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1423.963
$ws.Range("I19").Value = 845.5833
$ws.Range("J19").Value = 1886.6666
$ws.Range("K19").Value = 845.5833
$ws.Range("L19").Value = 1886.6666
$ws.Range("M19").Value = -670.5833
$ws.Range("N19").Value = -2236.6666

$ws.Range("H48").Value = 1892.5

$ws.Range("H56").Value = 1892.5

$ws.Range("H107").Value = 1420.5
$ws.Range("I107").Value = 1370.8422
$ws.Range("J107").Value = 1609.2
$ws.Range("K107").Value = 1370.8422
$ws.Range("L107").Value = 1609.2
$ws.Range("M107").Value = 549.1578
$ws.Range("N107").Value = -5449.2

$ws.Range("H112").Value = 2682541.2
$ws.Range("I112").Value = 1524
$ws.Range("J112").Value = 3486846.5
$ws.Range("K112").Value = 4572
$ws.Range("L112").Value = 10460539.5
$ws.Range("N112").Value = -10462755.5
$ws.Range("M112").Value = -3464

$ws.Range("H135").Value = 1689.6522
$ws.Range("I135").Value = 1588.5883
$ws.Range("K135").Value = 14297.2947
$ws.Range("M135").Value = -11762.2947

$ws.Range("H137").Value = 32340.553
$ws.Range("I137").Value = 63872.11
$ws.Range("J137").Value = 3962.15
$ws.Range("K137").Value = 191616.33
$ws.Range("L137").Value = 11886.45
$ws.Range("M137").Value = -189066.33
$ws.Range("N137").Value = -16986.45

$ws.Range("H138").Value = 2455.3572
$ws.Range("I138").Value = 1597.2222
$ws.Range("J138").Value = 4000
$ws.Range("K138").Value = 4791.6666
$ws.Range("L138").Value = 12000
$ws.Range("M138").Value = 348.3334000000004
$ws.Range("N138").Value = -22280

$ws.Range("H141").Value = 2217.6191
$ws.Range("I141").Value = 2318.4211
$ws.Range("J141").Value = 1260
$ws.Range("K141").Value = 6955.263300000001
$ws.Range("L141").Value = 3780
$ws.Range("M141").Value = -1775.263300000001
$ws.Range("N141").Value = -14140

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6053.154
$ws.Range("I32").Value = 2459.5652
$ws.Range("K32").Value = 2459.5652
$ws.Range("M32").Value = -2172.5652

$ws.Range("H74").Value = 42822.5
$ws.Range("I74").Value = 46529.09
$ws.Range("J74").Value = 2050
$ws.Range("K74").Value = 46529.09
$ws.Range("L74").Value = 2050
$ws.Range("M74").Value = -45655.09
$ws.Range("N74").Value = -3798

$ws.Range("H77").Value = 42822.5
$ws.Range("I77").Value = 46529.09
$ws.Range("J77").Value = 2050
$ws.Range("K77").Value = 232645.45
$ws.Range("L77").Value = 10250
$ws.Range("M77").Value = -228277.45
$ws.Range("N77").Value = -18986

$ws.Range("H97").Value = 1669.3549
$ws.Range("I97").Value = 1353.091
$ws.Range("J97").Value = 2442.4443
$ws.Range("K97").Value = 1353.091
$ws.Range("L97").Value = 2442.4443
$ws.Range("M97").Value = -857.0909999999999
$ws.Range("N97").Value = -3434.4443

$ws.Range("H110").Value = 4906.684
$ws.Range("I110").Value = 4758.6
$ws.Range("J110").Value = 5462
$ws.Range("K110").Value = 4758.6
$ws.Range("L110").Value = 5462
$ws.Range("M110").Value = -2713.6
$ws.Range("N110").Value = -9552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1790.4546
$ws.Range("I94").Value = 1115
$ws.Range("J94").Value = 2972.5
$ws.Range("K94").Value = 1115
$ws.Range("L94").Value = 2972.5
$ws.Range("M94").Value = -664
$ws.Range("N94").Value = -3874.5

$ws.Range("H132").Value = 123000
$ws.Range("J132").Value = 123000
$ws.Range("L132").Value = 123000
$ws.Range("N132").Value = -133120

$ws.Range("H134").Value = 2897.1973
$ws.Range("I134").Value = 3066.02
$ws.Range("K134").Value = 9198.059999999999
$ws.Range("M134").Value = -6663.059999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5858.7144
$ws.Range("I16").Value = 3999.3333
$ws.Range("J16").Value = 7253.25
$ws.Range("K16").Value = 3999.3333
$ws.Range("L16").Value = 7253.25
$ws.Range("M16").Value = -3712.3333
$ws.Range("N16").Value = -7827.25

$ws.Range("H31").Value = 129194.8
$ws.Range("I31").Value = 210480.52
$ws.Range("J31").Value = 3333.0322
$ws.Range("K31").Value = 210480.52
$ws.Range("L31").Value = 3333.0322
$ws.Range("M31").Value = -210185.52
$ws.Range("N31").Value = -3923.0322

$ws.Range("H34").Value = 129194.8
$ws.Range("I34").Value = 210480.52
$ws.Range("J34").Value = 3333.0322
$ws.Range("K34").Value = 210480.52
$ws.Range("L34").Value = 3333.0322
$ws.Range("M34").Value = -210278.52
$ws.Range("N34").Value = -3737.0322

$ws.Range("H58").Value = 3912
$ws.Range("I58").Value = 3125
$ws.Range("J58").Value = 5223.6665
$ws.Range("K58").Value = 3125
$ws.Range("L58").Value = 5223.6665
$ws.Range("M58").Value = -2922
$ws.Range("N58").Value = -5629.6665

$ws.Range("H99").Value = 285855
$ws.Range("I99").Value = 480887.9
$ws.Range("J99").Value = 12808.934
$ws.Range("K99").Value = 480887.9
$ws.Range("L99").Value = 12808.934
$ws.Range("M99").Value = -479389.9
$ws.Range("N99").Value = -15804.934

$ws.Range("H113").Value = 5858.7144
$ws.Range("I113").Value = 3999.3333
$ws.Range("J113").Value = 7253.25
$ws.Range("K113").Value = 3999.3333
$ws.Range("L113").Value = 7253.25
$ws.Range("M113").Value = -1829.3333
$ws.Range("N113").Value = -11593.25

$ws.Range("H126").Value = 285855
$ws.Range("I126").Value = 480887.9
$ws.Range("J126").Value = 12808.934
$ws.Range("K126").Value = 1442663.7
$ws.Range("L126").Value = 38426.802
$ws.Range("M126").Value = -1440193.7
$ws.Range("N126").Value = -43366.802

$ws.Range("H134").Value = 2166.4
$ws.Range("I134").Value = 2193.4614
$ws.Range("K134").Value = 6580.3842
$ws.Range("M134").Value = -4045.3842

$ws.Range("H136").Value = 3912
$ws.Range("I136").Value = 3125
$ws.Range("J136").Value = 5223.6665
$ws.Range("K136").Value = 9375
$ws.Range("L136").Value = 15670.9995
$ws.Range("M136").Value = -6825
$ws.Range("N136").Value = -20770.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 178.66667
$ws.Range("I23").Value = 160.66667
$ws.Range("J23").Value = 196.66667
$ws.Range("K23").Value = 482.00001
$ws.Range("L23").Value = 590.00001
$ws.Range("N23").Value = -1060.00001
$ws.Range("M23").Value = -247.00001

$ws.Range("H68").Value = 2877601.5
$ws.Range("I68").Value = 13891297
$ws.Range("J68").Value = 4463.478
$ws.Range("K68").Value = 41673891
$ws.Range("L68").Value = 13390.434
$ws.Range("M68").Value = -41673080
$ws.Range("N68").Value = -15012.434

$ws.Range("H71").Value = 2877601.5
$ws.Range("I71").Value = 13891297
$ws.Range("J71").Value = 4463.478
$ws.Range("K71").Value = 125021673
$ws.Range("L71").Value = 40171.302
$ws.Range("M71").Value = -125017617
$ws.Range("N71").Value = -48283.302

$ws.Range("H107").Value = 534.8461
$ws.Range("I107").Value = 496.08334
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 1488.25002
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 431.7499800000001
$ws.Range("N107").Value = -6840

$ws.Range("H137").Value = 1321.3334
$ws.Range("I137").Value = 1321.3334
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 3964.0002
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 1135.9998
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3774.738
$ws.Range("I80").Value = 3280.389
$ws.Range("J80").Value = 4145.5
$ws.Range("K80").Value = 3280.389
$ws.Range("L80").Value = 4145.5
$ws.Range("M80").Value = -2282.389
$ws.Range("N80").Value = -6141.5

$ws.Range("H83").Value = 3774.738
$ws.Range("I83").Value = 3280.389
$ws.Range("J83").Value = 4145.5
$ws.Range("K83").Value = 16401.945
$ws.Range("L83").Value = 20727.5
$ws.Range("M83").Value = -11409.945
$ws.Range("N83").Value = -30711.5

$ws.Range("H102").Value = 42036.69
$ws.Range("I102").Value = 3434.5
$ws.Range("J102").Value = 254348.75
$ws.Range("K102").Value = 3434.5
$ws.Range("L102").Value = 254348.75
$ws.Range("M102").Value = -1812.5
$ws.Range("N102").Value = -257592.75

$ws.Range("H126").Value = 20867.916
$ws.Range("I126").Value = 21854.637
$ws.Range("J126").Value = 10014
$ws.Range("K126").Value = 65563.91099999999
$ws.Range("L126").Value = 30042
$ws.Range("M126").Value = -63093.91099999999
$ws.Range("N126").Value = -34982

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5746.1523
$ws.Range("I40").Value = 5653.3057
$ws.Range("J40").Value = 6080.4
$ws.Range("K40").Value = 5653.3057
$ws.Range("L40").Value = 6080.4
$ws.Range("M40").Value = -5517.3057
$ws.Range("N40").Value = -6352.4

$ws.Range("H93").Value = 45456644
$ws.Range("I93").Value = 1113.6471
$ws.Range("J93").Value = 200005440
$ws.Range("K93").Value = 1113.6471
$ws.Range("L93").Value = 200005440
$ws.Range("M93").Value = 134.3529000000001
$ws.Range("N93").Value = -200007936

$ws.Range("H132").Value = 3580.7144
$ws.Range("I132").Value = 3641.4211
$ws.Range("K132").Value = 10924.2633
$ws.Range("M132").Value = -8394.263300000001

$ws.Range("H136").Value = 3797.8572
$ws.Range("I136").Value = 3285.6
$ws.Range("J136").Value = 8066.6665
$ws.Range("K136").Value = 9856.799999999999
$ws.Range("L136").Value = 24199.9995
$ws.Range("M136").Value = -7306.799999999999
$ws.Range("N136").Value = -29299.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1184.4375
$ws.Range("I113").Value = 922.36365
$ws.Range("K113").Value = 2767.09095
$ws.Range("M113").Value = -597.0909499999998

$ws.Range("H126").Value = 2578.0952
$ws.Range("I126").Value = 2361.353
$ws.Range("J126").Value = 3499.25
$ws.Range("K126").Value = 7084.059
$ws.Range("L126").Value = 10497.75
$ws.Range("M126").Value = -4614.059
$ws.Range("N126").Value = -15437.75

$ws.Range("H132").Value = 1979.22
$ws.Range("I132").Value = 999.6842
$ws.Range("J132").Value = 2579.5806
$ws.Range("K132").Value = 2999.0526
$ws.Range("L132").Value = 7738.7418
$ws.Range("M132").Value = -469.0526
$ws.Range("N132").Value = -12798.7418

$ws.Range("H136").Value = 358420.72
$ws.Range("I136").Value = 401071.2
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 1203213.6
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -1200663.6
$ws.Range("N136").Value = -14100

Write-Output "Applied 267 value updates and 1 cell clears."
